$d = $word.ActiveDocument

# The document's final paragraph (before the sectPr) currently holds only a
# tab stop definition and a single tab run. We need to:
#   1. Insert a new Heading1 paragraph "Metal Switch Body" right before it.
#   2. Prepend two new runs of body text to that final paragraph (ahead of
#      its existing tab run) and drop its custom tab-stop paragraph formatting.

$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Last
$lastRange = $lastPara.Range

# 1) Insert the new heading paragraph before the last paragraph.
$lastRange.InsertParagraphBefore()

# Re-fetch paragraphs after the structural change.
$headingPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$headingRange = $headingPara.Range
$headingRange.Text = "Metal Switch Body"
$headingRange.ParagraphFormat.Style = "Heading1"

# 2) Prepend the body-text runs to the true last paragraph, ahead of its
#    existing tab character, and clear its tab-stop formatting.
$finalPara = $d.Paragraphs.Last
$finalRange = $finalPara.Range
$insertPoint = $d.Range($finalRange.Start, $finalRange.Start)
$insertPoint.InsertBefore("It was found that for some use cases, the plastic switch body used in the original switch could shatter with heavy use, and some users wanted the cable to come out the same side as the hinge. A metal switch body was sourced on Amazon and added to the BoM and Maker Guide as an alternate design.")

$finalPara.Range.ParagraphFormat.TabStops.ClearAll()
